# Updates the cryptos price/volume table (columns D and E) on the active sheet
# to reflect the latest refreshed values from the GitHub Actions data pull.
# Values that look like plain numbers (e.g. "605.53") are entered with a
# leading apostrophe so Excel keeps them as text, matching how this sheet
# already stores its Price column (it mixes thousand-dot formatted figures
# like "66.807.88" with plain decimals like "605.53", all stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.807.88"
$ws.Range("E2").Value = "  +2.88%  "

$ws.Range("D3").Value = "3.199.84"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'605.53"
$ws.Range("E5").Value = "  +4.55%  "

$ws.Range("D6").Value = "'156.98"
$ws.Range("E6").Value = "  +5.67%  "

$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = "  +6.10%  "

$ws.Range("D9").Value = "3.198.01"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("D11").Value = "'5.93"
$ws.Range("E11").Value = "  -3.37%  "

$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "'39.22"
$ws.Range("E14").Value = "  +5.69%  "

$ws.Range("D15").Value = "3.724.36"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").Value = "66.745.61"
$ws.Range("E16").Value = "  +2.82%  "

$ws.Range("D17").Value = "'7.48"
$ws.Range("E17").Value = "  +4.82%  "

$ws.Range("D18").Value = "3.201.22"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("D20").Value = "'522.09"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").Value = "'15.45"
$ws.Range("E21").Value = "  +2.42%  "

$ws.Range("E22").Value = "  +4.17%  "

$ws.Range("E23").Value = "  +6.31%  "

$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'85.44"
$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("E27").Value = "  +2.01%  "

$ws.Range("D28").Value = "'3.03"
$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("E29").Value = "  +9.58%  "

$ws.Range("E30").Value = "  +7.78%  "

$ws.Range("D31").Value = "'7.03"
$ws.Range("E31").Value = "  +9.16%  "

$ws.Range("D32").Value = "'28.36"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("E33").Value = "  +3.33%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "'6.57"
$ws.Range("E35").Value = "  +1.50%  "

$ws.Range("D36").Value = "'520.33"
$ws.Range("E36").Value = "  +9.55%  "

$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("D39").Value = "'0.0426"
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("E40").Value = "  +8.64%  "

$ws.Range("E41").Value = "  +2.04%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("E43").Value = "  +15.75%  "

$ws.Range("E44").Value = "  +7.26%  "

$ws.Range("E45").Value = "  +2.09%  "

$ws.Range("D46").Value = "2.897.46"
$ws.Range("E46").Value = "  -3.50%  "

$ws.Range("E47").Value = "  +1.53%  "

$ws.Range("E48").Value = "  +11.22%  "

$ws.Range("E49").Value = "  +3.86%  "

$ws.Range("E51").Value = "  +4.44%  "
